# Generate Report for Handoff
# - Refreshes the "current" handoff file (old UUID -> new UUID, new dates)
# - Adds a brand-new handback row (new UUID) to the Overview / zh-cn / de-de sheets

$wb = $excel.ActiveWorkbook

$oldGuid = "f24e5f39-3124-44b0-857f-12984381c4d2"
$newGuid = "2477acd5-bbfe-449a-a14c-cf5a121f0e21"
$addGuid = "bc81d3b0-5149-4471-80bb-638b501daaaa"

$newGuidHash = "30e28e2e0e3301109b173c0d68b0e2f3f7f30832"
$addGuidHash = "c62a57a6d99adc1becb85007cdbd5ae5f74b279d"

$oldDateOverview = "2017-02-21 05:09:51"
$newDateOverview = "2017-02-21 05:10:43"

$oldDateZh = "2017-02-21 05:09:37"
$newDateZh = "2017-02-21 05:10:27"

$oldDateDe = "2017-02-21 05:09:51"
$newDateDe = "2017-02-21 05:10:43"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/5a362cc0871c724b405ed115d5ab33e4d0917358/e2e/"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$lo = $ws.ListObjects.Item(1)

# update row 2 (existing handoff entry)
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("B2").Value = "e2e\$newGuid.md"
$ws.Range("B2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "$repoBase$newGuid.md", "", "", "e2e\$newGuid.md") | Out-Null
$ws.Range("G2").Value = $newDateOverview

# grow the table and append the new row
$lo.ListRows.Add() | Out-Null
$ws.Range("A3").Value = "$addGuid.md"
$ws.Range("B3").Value = "e2e\$addGuid.md"
$ws.Hyperlinks.Add($ws.Range("B3"), "$repoBase$addGuid.md", "", "", "e2e\$addGuid.md") | Out-Null
$ws.Range("C3").Value = ".md"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = $newDateOverview

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$lo = $ws.ListObjects.Item(1)

# update row 2 (existing handoff entry)
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase$newGuid.md", "", "", "$newGuid.md") | Out-Null
$ws.Range("G2").Value = "$newGuid.$newGuidHash.zh-cn.xlf"
$ws.Range("H2").Value = $newDateZh

# grow the table and append the new row
$lo.ListRows.Add() | Out-Null
$ws.Range("A3").Value = "$addGuid.md"
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase$addGuid.md", "", "", "$addGuid.md") | Out-Null
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "$addGuid.$addGuidHash.zh-cn.xlf"
$ws.Range("H3").Value = $newDateZh
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "0001-01-01 00:00:00"
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "True"
$ws.Range("P3").Value = ""
$ws.Range("Q3").Value = "False"
$ws.Range("R3").Value = ""

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$lo = $ws.ListObjects.Item(1)

# update row 2 (existing handoff entry)
$ws.Range("A2").Value = "$newGuid.md"
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "$repoBase$newGuid.md", "", "", "$newGuid.md") | Out-Null
$ws.Range("G2").Value = "$newGuid.$newGuidHash.de-de.xlf"
$ws.Range("H2").Value = $newDateDe

# grow the table and append the new row
$lo.ListRows.Add() | Out-Null
$ws.Range("A3").Value = "$addGuid.md"
$ws.Hyperlinks.Add($ws.Range("A3"), "$repoBase$addGuid.md", "", "", "$addGuid.md") | Out-Null
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "e2e"
$ws.Range("E3").Value = "ht"
$ws.Range("F3").Value = "False"
$ws.Range("G3").Value = "$addGuid.$addGuidHash.de-de.xlf"
$ws.Range("H3").Value = $newDateDe
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("L3").Value = "0001-01-01 00:00:00"
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("O3").Value = "True"
$ws.Range("P3").Value = ""
$ws.Range("Q3").Value = "False"
$ws.Range("R3").Value = ""

Write-Host "Generate Report for Handoff: done"
